# Applies the "vex / smart1 / giotto" namespace-registry update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "Last Update" counter / date in row 1 (B1: 2024-11-22 -> 2025-02-27)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = 45715

# ---------------------------------------------------------------------------
# 2. Row 14 ("ama" entry) loses its special one-off formatting and is
#    restyled to match the plain table rows (e.g. row 15).
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).ClearFormats()
$ws.Range("A15:U15").Copy()
$ws.Range("A14:U14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Insert three new rows right before the old row 127 (the blank separator
#    above the "Planetary missions" section) for the new PSA mission
#    namespaces: vex, smart1, giotto. Everything below shifts down by 3.
# ---------------------------------------------------------------------------
$ws.Rows.Item(127).Insert()
$ws.Rows.Item(128).Insert()
$ws.Rows.Item(129).Insert()

# Copy the formatting of the most similar existing PSA row (126) onto the
# three freshly inserted rows.
$ws.Range("A126:U126").Copy()
$ws.Range("A127:U129").PasteSpecial(-4122)

# --- Row 127: vex ----------------------------------------------------------
$ws.Range("A127").Value = "vex"
$ws.Range("B127").Value = "Venus Express Mission"
$ws.Range("C127").Value = "Namespace for the Venus Express Mission schema."
$ws.Range("D127").Value = "vex"
$ws.Range("E127").Value = "http://psa.esa.int/psa/vex/v1"
$ws.Range("F127").Value = "vex"
$ws.Range("G127").Value = "urn:esa:psa"
$ws.Range("H127").Value = "PDS4_PSA_VEX"
$ws.Range("I127").Value = "Mission"
$ws.Range("J127").Value = "0001_ESA_PSA_1"
$ws.Range("K127").Value = "vex"
$ws.Range("L127").Value = "vex"
$ws.Range("M127").Value = "vex"
$ws.Range("N127").Value = "Mark Bentley"
$ws.Range("O127").Value = "Mark.Bentley at esa.int"
$ws.Range("Q127").Value = 45715
$ws.Range("R127").Value = "M. Bentley"
$ws.Range("S127").Value = "No"
$ws.Range("T127").Value = "No"

# --- Row 128: smart1 --------------------------------------------------------
$ws.Range("A128").Value = "smart1"
$ws.Range("B128").Value = "Small Missions for Advanced Research in Technology Mission"
$ws.Range("C128").Value = "Namespace for the Small Missions for Advanced Research in Technology mission schema."
$ws.Range("D128").Value = "smart1"
$ws.Range("E128").Value = "http://psa.esa.int/psa/smart1/v1"
$ws.Range("F128").Value = "smart1"
$ws.Range("G128").Value = "urn:esa:psa"
$ws.Range("H128").Value = "PDS4_PSA_SMART1"
$ws.Range("I128").Value = "Mission"
$ws.Range("J128").Value = "0001_ESA_PSA_1"
$ws.Range("K128").Value = "smart1"
$ws.Range("L128").Value = "smart1"
$ws.Range("M128").Value = "smart1"
$ws.Range("N128").Value = "Mark Bentley"
$ws.Range("O128").Value = "Mark.Bentley at esa.int"
$ws.Range("Q128").Value = 45715
$ws.Range("R128").Value = "M. Bentley"
$ws.Range("S128").Value = "No"
$ws.Range("T128").Value = "No"
$ws.Rows.Item(128).RowHeight = 48

# --- Row 129: giotto ---------------------------------------------------------
$ws.Range("A129").Value = "giotto"
$ws.Range("B129").Value = "Giotto Mission"
$ws.Range("C129").Value = "Namespace for the Giotto Mission schema."
$ws.Range("D129").Value = "giotto"
$ws.Range("E129").Value = "http://psa.esa.int/psa/giotto/v1"
$ws.Range("F129").Value = "giotto"
$ws.Range("G129").Value = "urn:esa:psa"
$ws.Range("H129").Value = "PDS4_PSA_GIOTTO"
$ws.Range("I129").Value = "Mission"
$ws.Range("J129").Value = "0001_ESA_PSA_1"
$ws.Range("K129").Value = "giotto"
$ws.Range("L129").Value = "giotto"
$ws.Range("M129").Value = "giotto"
$ws.Range("N129").Value = "Mark Bentley"
$ws.Range("O129").Value = "Mark.Bentley at esa.int"
$ws.Range("Q129").Value = 45715
$ws.Range("R129").Value = "M. Bentley"
$ws.Range("S129").Value = "No"
$ws.Range("T129").Value = "No"
